$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the "Sentence Structure" (B) and "Action" (D) values      ---
# --- (together with B's formatting) between rows 17 (TR16) and 18 (TR17).   ---
# Stage B17's original formatting in a scratch cell so it can be moved onto
# B18 after B18's formatting has been copied onto B17.
$ws.Range("B17").Copy()
$ws.Range("AA1").PasteSpecial(-4122)

$ws.Range("B18").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("AA1").Copy()
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("B17").Value = "AO"
$ws.Range("B18").Value = "SVCVO"
$ws.Range("D17").Value = "actor=A, target=C, action=B"
$ws.Range("D18").Value = "actor=B, target=D, action=C"

$ws.Range("AA1").Clear()

# --- Step 2: append new transformation rules TR18-TR22 as rows 19-23, ---
# --- reusing row 18's current (post-swap) formatting as the template. ---
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)
$ws.Range("A20:D20").PasteSpecial(-4122)
$ws.Range("A21:D21").PasteSpecial(-4122)
$ws.Range("A22:D22").PasteSpecial(-4122)
$ws.Range("A23:D23").PasteSpecial(-4122)

$ws.Range("A19").Value = "TR18"
$ws.Range("B19").Value = "SVVO"
$ws.Range("C19").Value = "add_behavior"
$ws.Range("D19").Value = "actor=B, target=D, action=C"

$ws.Range("A20").Value = "TR19"
$ws.Range("B20").Value = "SVVOcomma1"
$ws.Range("C20").Value = "add_behavior"
$ws.Range("D20").Value = "actor=B, target=D, action=C"

$ws.Range("A21").Value = "TR20"
$ws.Range("B21").Value = "SVOCO"
$ws.Range("C21").Value = "add_behavior"
$ws.Range("D21").Value = "actor=B, target=CD, action=A"

$ws.Range("A22").Value = "TR21"
$ws.Range("B22").Value = "SVOO"
$ws.Range("C22").Value = "add_behavior"
$ws.Range("D22").Value = "actor=B, target=CD, action=A"

$ws.Range("A23").Value = "TR22"
$ws.Range("B23").Value = "SVOOcomma1"
$ws.Range("C23").Value = "add_behavior"
$ws.Range("D23").Value = "actor=B, target=CD, action=A"

# --- Step 3: B23 uses a distinct font (Arial 10, black) instead of the ---
# --- usual Arial 10 theme-color font used by the other new rows.      ---
$ws.Range("B6").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Font.Color = 0
$ws.Range("B23").Value = "SVOOcomma1"

# --- Step 4: restore the active selection the author ended up with. ---
$null = $ws.Range("C26").Select()
